$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.491.20"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.553.54"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.483"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.28"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0582"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0892"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.777.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.549.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.503.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.510"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.73"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.22"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("E30").Value = "  -4.69%  "
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.387.68"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.65"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.516"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0463"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.30"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.690.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +2.45%  "
